$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source added two new price records right after the existing row 463,
# pushing the old rows 464:574 down to 466:576 (dimension grows from R574 to R576).
$ws.Range("A464:R465").EntireRow.Insert()

# New row 464
$ws.Range("A464").Value = 6
$ws.Range("B464").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C464").Value = "Metropolitana"
$ws.Range("D464").Value = 44508
$ws.Range("E464").Value = 13
$ws.Range("F464").Value = 100112040
$ws.Range("G464").Value = "Cilantro"
$ws.Range("H464").Value = "Sin especificar"
$ws.Range("I464").Value = "Primera"
$ws.Range("J464").Value = 430
$ws.Range("K464").Value = 4500
$ws.Range("L464").Value = 5000
$ws.Range("M464").Value = 4709
$ws.Range("N464").Value = "`$/caja 36 atados"
$ws.Range("O464").Value = "Región Metropolitana"
$ws.Range("P464").Value = 131
$ws.Range("Q464").Value = 36
$ws.Range("R464").Value = "Hortaliza"

# New row 465
$ws.Range("A465").Value = 6
$ws.Range("B465").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C465").Value = "Metropolitana"
$ws.Range("D465").Value = 44508
$ws.Range("E465").Value = 13
$ws.Range("F465").Value = 100112040
$ws.Range("G465").Value = "Cilantro"
$ws.Range("H465").Value = "Sin especificar"
$ws.Range("I465").Value = "Primera"
$ws.Range("J465").Value = 330
$ws.Range("K465").Value = 9000
$ws.Range("L465").Value = 10000
$ws.Range("M465").Value = 9455
$ws.Range("N465").Value = "`$/docena de atados"
$ws.Range("O465").Value = "Región Metropolitana"
$ws.Range("P465").Value = 3152
$ws.Range("Q465").Value = 3
$ws.Range("R465").Value = "Hortaliza"
